# Adding Screws for Servos
# Insert a new row above row 6 (shifts rows 6-14 -> 7-15) and fill it in
# with the new "Servo Screws" line item, then fix up the hyperlinks that
# Insert() leaves anchored to their old (pre-shift) cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("6:6").Insert()

# Match the formatting used by the rest of the data rows (font Arial 10,
# same as A7:B14 etc.) for the new item/price cells.
$ws.Range("A6:B6").Font.Name = "Arial"
$ws.Range("A6:B6").Font.Size = 10

$ws.Range("A6").Value = "Servo Screws"
$ws.Range("B6").Value = 5.9

# Insert() does not move the existing hyperlinks along with the cells they
# were anchored to, so rebuild the whole hyperlink list with the correct,
# post-shift cell references (D6 old contents now live on D7, etc.), plus
# the brand new hyperlink for D6 ("Servo Screws" -> McMaster link).
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("D2"), "http://www.robotis.us/opencm9-04-c-with-onboard-xl-type-connectors/")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://www.amazon.com/MJX-F645-Replacement-Battery-Meijiaxin/dp/B007XOMXMO/ref=sr_1_35?ie=UTF8&qid=1526604820&sr=8-35&keywords=7.4v+1500mah")
$ws.Hyperlinks.Add($ws.Range("D4"), "https://www.amazon.com/Performix-11602-6-075815116024-Yellow-Plasti/dp/B000ZN1T16/ref=sr_1_13?ie=UTF8&qid=1526518098&sr=8-13&keywords=plastidip")
$ws.Hyperlinks.Add($ws.Range("D7"), "https://www.mcmaster.com/", "92005a076/=1cvfwje")
$ws.Hyperlinks.Add($ws.Range("D8"), "http://www.robotis.us/dynamixel-xl-320/")
$ws.Hyperlinks.Add($ws.Range("D9"), "https://www.amazon.com/QWinOut-Switching-Helicopter-Quadcopter-Hexacopter/dp/B01D10MWYW/ref=sr_1_3?s=toys-and-games&ie=UTF8&qid=1526605943&sr=1-3&keywords=liion+battery+charger&refinements=p_36%3A1253560011")
$ws.Hyperlinks.Add($ws.Range("D10"), "https://www.amazon.com/eBoot-Connector-Female-Cable-Battery/dp/B01M5AHF0Z")
$ws.Hyperlinks.Add($ws.Range("D11"), "https://www.mcmaster.com/", "acrylic/=1cvfzgg")
$ws.Hyperlinks.Add($ws.Range("D12"), "https://www.mcmaster.com/", "93625a102/=1cvfxju")
$ws.Hyperlinks.Add($ws.Range("D13"), "http://www.robotis.us/bt-210/")
$ws.Hyperlinks.Add($ws.Range("D14"), "https://www.mcmaster.com/", "93657a212/=1cvg4m2")
$ws.Hyperlinks.Add($ws.Range("D6"), "https://www.mcmaster.com/92005a078")

# Hyperlinks.Add() stamps a fresh "hyperlink-like" style instead of reusing
# the sheet's existing Hyperlink cell style, so reapply it explicitly on
# every linked cell to keep the style sheet / cell formatting unchanged.
$ws.Range("D2").Style = "Hyperlink"
$ws.Range("D3").Style = "Hyperlink"
$ws.Range("D4").Style = "Hyperlink"
$ws.Range("D6").Style = "Hyperlink"
$ws.Range("D7").Style = "Hyperlink"
$ws.Range("D8").Style = "Hyperlink"
$ws.Range("D9").Style = "Hyperlink"
$ws.Range("D10").Style = "Hyperlink"
$ws.Range("D11").Style = "Hyperlink"
$ws.Range("D12").Style = "Hyperlink"
$ws.Range("D13").Style = "Hyperlink"
$ws.Range("D14").Style = "Hyperlink"

$ws.Range("D6").Select()
